# The "IVY ZAD SYRUP 120 ML" line item (row 57) was removed from the Day
# Sale report. Deleting the entire row shifts every row below it up by
# one, which Excel handles natively (including re-numbering the shared
# string table and the merged cell ranges), matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(57).Delete()

# The grand-total cell (now on row 129 after the shift) was a hard-coded
# value, not a formula, so it still reflects the old total. Subtract the
# deleted row's sale amount (qty 1 x 60.00 = 60.00) to match the new
# total reported in the refreshed export.
$ws.Range("P129").Value = 7149.6499999999996

# The footer (now on row 130) carries the report's "generated at"
# timestamp; the new upload was produced four minutes later.
$ws.Range("A130").Value = "Wednesday, 16 July, 2025 10:17 PM"
